# Casos de prueba.xlsx - apply corrections to the test-case table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Content fixes -------------------------------------------------------
# Caso de prueba 1: el valor de entrada cambia de "1,2,3" a "7,5,3"
$ws.Range("B2").Value = "7,5,3"

# Casos de prueba 9 y 10: se quita el punto final del mensaje obtenido
$ws.Range("D10").Value = "La suma de los primeros dos lados es menor al tercero"
$ws.Range("D11").Value = "La suma de los primeros dos lados es menor al tercero"

# --- Formato: marcar que casos pasaron (verde) vs los que faltan (naranja) -
# Todas las filas de datos de la tabla de casos de prueba (filas 2 a 14, col B:D)
$green = 0 + (176 * 256) + (80 * 65536)
$orange = 255 + (192 * 256) + (0 * 65536)

# Casos que ya funcionan -> verde (filas 2-8 y 12-14)
$ws.Range("B2:D8").Interior.Color = $green
$ws.Range("B12:D14").Interior.Color = $green

# Casos de prueba 8, 9 y 10 (filas 9, 10, 11) siguen sin funcionar -> naranja
$failRange = $ws.Range("B9:D11")
$failRange.Interior.Color = $orange

# --- Seleccion activa ------------------------------------------------------
$ws.Range("B3").Select()
